$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.715.77'
$ws.Range("E2").Value = '  +0.38%  '
$ws.Range("D3").Value = '1.600.56'
$ws.Range("E3").Value = '  +0.24%  '
$ws.Range("E4").Value = '  +0.34%  '
$__style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.47'
$ws.Range("D5").Style = $__style
$ws.Range("E7").Value = '  +0.33%  '
$ws.Range("E8").Value = '  +0.10%  '
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("E10").Value = '  +0.66%  '
$__style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0843'
$ws.Range("D11").Style = $__style
$ws.Range("E11").Value = '  +0.33%  '
$ws.Range("D12").Value = '1.825.69'
$ws.Range("E12").Value = '  +0.29%  '
$ws.Range("D13").Value = '1.591.22'
$ws.Range("E13").Value = '  -0.70%  '
$__style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.524'
$ws.Range("D15").Style = $__style
$ws.Range("E15").Value = '  +0.30%  '
$__style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.00'
$ws.Range("D16").Style = $__style
$ws.Range("E16").Value = '  +0.19%  '
$ws.Range("D17").Value = '26.685.23'
$ws.Range("E17").Value = '  +0.31%  '
$ws.Range("D18").Value = '0.0₃0745'
$ws.Range("E18").Value = '  +0.79%  '
$__style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '210.29'
$ws.Range("D19").Style = $__style
$ws.Range("E19").Value = '  +1.02%  '
$ws.Range("E20").Value = '  +2.77%  '
$ws.Range("E21").Value = '  +0.32%  '
$ws.Range("E22").Value = '  +0.10%  '
$ws.Range("E23").Value = '  -0.66%  '
$__style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.23'
$ws.Range("D25").Style = $__style
$ws.Range("E26").Value = '  +0.25%  '
$ws.Range("E27").Value = '  +0.08%  '
$ws.Range("E28").Value = '  -0.97%  '
$__style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.37'
$ws.Range("D29").Style = $__style
$ws.Range("E29").Value = '  +0.62%  '
$__style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0512'
$ws.Range("D30").Style = $__style
$ws.Range("E30").Value = '  +0.24%  '
$ws.Range("E31").Value = '  +0.07%  '
$ws.Range("E32").Value = '  +1.21%  '
$ws.Range("E33").Value = '  +0.72%  '
$ws.Range("D34").Value = '1.296.69'
$ws.Range("E34").Value = '  +0.98%  '
$__style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.48'
$ws.Range("D35").Style = $__style
$ws.Range("E35").Value = '  +1.00%  '
$__style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.606'
$ws.Range("D36").Style = $__style
$ws.Range("E36").Value = '  -2.15%  '
$ws.Range("E37").Value = '  +0.64%  '
$ws.Range("E38").Value = '  +14.64%  '
$__style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.825'
$ws.Range("D40").Style = $__style
$ws.Range("E40").Value = '  -1.76%  '
$__style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.40'
$ws.Range("D41").Style = $__style
$ws.Range("E41").Value = '  -1.31%  '
$ws.Range("E42").Value = '  -0.09%  '
$ws.Range("E43").Value = '  -0.70%  '
$__style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '63.00'
$ws.Range("D44").Style = $__style
$ws.Range("E44").Value = '  -1.97%  '
$ws.Range("D45").Value = '1.738.64'
$ws.Range("E45").Value = '  +0.36%  '
$__style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '90.60'
$ws.Range("D46").Style = $__style
$ws.Range("E46").Value = '  +0.55%  '
$__style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.56'
$ws.Range("D47").Style = $__style
$ws.Range("E47").Value = '  -2.71%  '
$ws.Range("E48").Value = '  -1.01%  '
$__style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0517'
$ws.Range("D49").Style = $__style
$ws.Range("E49").Value = '  +1.86%  '
$ws.Range("E50").Value = '  +0.04%  '
$ws.Range("E51").Value = '  -0.07%  '
